$d = $word.ActiveDocument

# 1. Delete the "Meta description: ..." paragraph near the top of the document
#    (paragraph 2, right after the title heading).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new bold paragraph "Play Buffalo Bounty Free: Game Review & Features"
#    right before the last paragraph (the italic "Prompt:" paragraph).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore() | Out-Null
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Buffalo Bounty Free: Game Review &amp; Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xml) | Out-Null

# 3. Replace the text of the last (italic) paragraph from the old "Prompt: ..." text
#    to the new meta-description sentence, keeping its italic formatting.
$old = "Prompt: Create a feature image for Buffalo Bounty. The image should fit the game's theme and be in a cartoon style. The image should feature a happy Maya warrior wearing glasses. This warrior should be in a majestic background surrounded by mountains, wildlife and maybe even a buffalo. The image should be colorful and eye-catching, capturing the essence of the game's rustic landscape and adventurous spirit."
$new = "Experience the beautiful wilderness of North America in Buffalo Bounty, an online slot game featuring free spins and multiple power-up spin options."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
